# Split the single "Holidays are a list, ... dummy variables. " run into
# three runs:
#   1) "Holidays are a list, provided by the user,"
#   2) " "                                                  (one space)
#   3) "which mark specific dates as holidays and are dummy variables. "
# (the original had a double space between "user," and "which"; after the
# edit only a single space remains, now living in its own run)
# All three runs keep the exact same run formatting (rPr) as the original.

$d = $word.ActiveDocument

$oldText = "Holidays are a list, provided by the user,  which mark specific dates as holidays and are dummy variables. "

$searchRange = $d.Content
$searchRange.Find.ClearFormatting()
$found = $searchRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "edit.ps1: target sentence was not found in the document"
}

# Re-seat into a plain Range (a Range that has been used as a Find/Execute
# target behaves as an insertion point for InsertXML; a fresh Range with the
# same Start/End correctly replaces the whole matched span instead).
$runRange = $d.Range($searchRange.Start, $searchRange.End)

$part1 = "Holidays are a list, provided by the user,"
$part2 = " "
$part3 = "which mark specific dates as holidays and are dummy variables. "

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>$rPr<w:t>$part1</w:t></w:r><w:r>$rPr<w:t xml:space="preserve">$part2</w:t></w:r><w:r>$rPr<w:t xml:space="preserve">$part3</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$runRange.InsertXML($xml)
